# Auto-generated Excel COM-interop script
# Applies targeted numeric-cell updates across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the authoritative diff of Kujata_Profits market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 33: Glazed and Confused | Clear Glass Lens
  $ws.Range("H33").Value = 435.8889
  $ws.Range("I33").Value = 413.2
  $ws.Range("K33").Value = 413.2
  $ws.Range("M33").Value = -184.2
  # Row 45: The House Always Wins | Blinding Potion
  $ws.Range("H45").Value = 0
  $ws.Range("J45").Value = 0
  $ws.Range("L45").Value = 0
  $ws.Range("N45").ClearContents()
  # Row 132: Fast-forwarding Flora | Growth Formula Lambda
  $ws.Range("H132").Value = 5752699.5
  $ws.Range("I132").Value = 7940376.5
  $ws.Range("J132").Value = 10048
  $ws.Range("K132").Value = 23821129.5
  $ws.Range("L132").Value = 30144
  $ws.Range("M132").Value = -23818599.5
  $ws.Range("N132").Value = -35204
  # Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
  $ws.Range("H137").Value = 1435.88
  $ws.Range("I137").Value = 878.9091
  $ws.Range("K137").Value = 2636.7273
  $ws.Range("M137").Value = -86.72730000000001
  # Row 138: All-night Crafting | Cunning Craftsman's Tisane
  $ws.Range("H138").Value = 464231.1
  $ws.Range("I138").Value = 1596.0454
  $ws.Range("J138").Value = 605591.8
  $ws.Range("K138").Value = 4788.1362
  $ws.Range("L138").Value = 1816775.4
  $ws.Range("M138").Value = 351.8638000000001
  $ws.Range("N138").Value = -1827055.4

$ws = $wb.Worksheets.Item("ARM")
  # Row 32: Ingot We Trust | Steel Ingot
  $ws.Range("H32").Value = 4973.1
  $ws.Range("I32").Value = 4447.6
  $ws.Range("J32").Value = 9702.6
  $ws.Range("K32").Value = 4447.6
  $ws.Range("L32").Value = 9702.6
  $ws.Range("M32").Value = -4160.6
  $ws.Range("N32").Value = -10276.6
  # Row 74: As the Bolt Flies | Titanium Nugget
  $ws.Range("H74").Value = 1686.16
  $ws.Range("I74").Value = 946.0769
  $ws.Range("K74").Value = 946.0769
  $ws.Range("M74").Value = -72.07690000000002
  # Row 77: Heavy Metal Banned (L) | Titanium Nugget
  $ws.Range("H77").Value = 1686.16
  $ws.Range("I77").Value = 946.0769
  $ws.Range("K77").Value = 4730.3845
  $ws.Range("M77").Value = -362.3845000000001
  # Row 122: Haste for High Durium | High Durium Nugget
  $ws.Range("H122").Value = 1368.2
  $ws.Range("I122").Value = 1241.3334
  $ws.Range("J122").Value = 1748.8
  $ws.Range("K122").Value = 3724.0002
  $ws.Range("L122").Value = 5246.4
  $ws.Range("M122").Value = -1274.0002
  $ws.Range("N122").Value = -10146.4
  # Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
  $ws.Range("H132").Value = 1603.5652
  $ws.Range("I132").Value = 1362.4
  $ws.Range("K132").Value = 4087.2
  $ws.Range("M132").Value = -1557.2

$ws = $wb.Worksheets.Item("BSM")
  # Row 82: Spirituality Inspector | Titanium Lump Hammer
  $ws.Range("H82").Value = 20314.25
  $ws.Range("J82").Value = 31000
  $ws.Range("L82").Value = 31000
  $ws.Range("N82").Value = -31766
  # Row 85: The Clamor for Hammers (L) | Titanium Lump Hammer
  $ws.Range("H85").Value = 20314.25
  $ws.Range("J85").Value = 31000
  $ws.Range("L85").Value = 31000
  $ws.Range("N85").Value = -33652
  # Row 86: Through Thick and Thin | Adamantite Nugget
  $ws.Range("H86").Value = 4211.9165
  $ws.Range("I86").Value = 4450.3335
  $ws.Range("J86").Value = 3496.6667
  $ws.Range("K86").Value = 4450.3335
  $ws.Range("L86").Value = 3496.6667
  $ws.Range("M86").Value = -3327.3335
  $ws.Range("N86").Value = -5742.6667
  # Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
  $ws.Range("H89").Value = 4211.9165
  $ws.Range("I89").Value = 4450.3335
  $ws.Range("J89").Value = 3496.6667
  $ws.Range("K89").Value = 22251.6675
  $ws.Range("L89").Value = 17483.3335
  $ws.Range("M89").Value = -16635.6675
  $ws.Range("N89").Value = -28715.3335

$ws = $wb.Worksheets.Item("CRP")
  # Row 2: In with the New | Bone Harpoon
  $ws.Range("H2").Value = 0
  $ws.Range("I2").Value = 0
  $ws.Range("K2").Value = 0
  $ws.Range("M2").ClearContents()
  # Row 31: Wall Not Found | Walnut Lumber
  $ws.Range("H31").Value = 2002
  $ws.Range("I31").Value = 2093.3333
  $ws.Range("J31").Value = 1962.8572
  $ws.Range("K31").Value = 2093.3333
  $ws.Range("L31").Value = 1962.8572
  $ws.Range("M31").Value = -1798.3333
  $ws.Range("N31").Value = -2552.8572
  # Row 34: Armoires of the Rich and Famous | Walnut Lumber
  $ws.Range("H34").Value = 2002
  $ws.Range("I34").Value = 2093.3333
  $ws.Range("J34").Value = 1962.8572
  $ws.Range("K34").Value = 2093.3333
  $ws.Range("L34").Value = 1962.8572
  $ws.Range("M34").Value = -1891.3333
  $ws.Range("N34").Value = -2366.8572
  # Row 56: Trident and Error | Cobalt Trident
  $ws.Range("H56").Value = 0
  $ws.Range("J56").Value = 0
  $ws.Range("L56").Value = 0
  $ws.Range("N56").ClearContents()
  # Row 132: Hull Lotta Damage | Ginseng Lumber
  $ws.Range("H132").Value = 1557.7084
  $ws.Range("I132").Value = 1205.2106
  $ws.Range("K132").Value = 3615.6318
  $ws.Range("M132").Value = -1085.6318

$ws = $wb.Worksheets.Item("CUL")
  # Row 33: Cooking with Gas | Chicken Stock
  $ws.Range("H33").Value = 292.57895
  $ws.Range("I33").Value = 214.54546
  $ws.Range("J33").Value = 399.875
  $ws.Range("K33").Value = 1287.27276
  $ws.Range("L33").Value = 2399.25
  $ws.Range("M33").Value = -1004.27276
  $ws.Range("N33").Value = -2965.25
  # Row 34: Fever Pitch | Chamomile Tea
  $ws.Range("H34").Value = 2423.875
  $ws.Range("I34").Value = 2133.3333
  $ws.Range("J34").Value = 2598.2
  $ws.Range("K34").Value = 6399.999899999999
  $ws.Range("L34").Value = 7794.599999999999
  $ws.Range("M34").Value = -6315.999899999999
  $ws.Range("N34").Value = -7962.599999999999
  # Row 39: Bloody Good Tart, This | Blood Currant Tart
  $ws.Range("H39").Value = 2493.1765
  $ws.Range("J39").Value = 2213.1428
  $ws.Range("L39").Value = 6639.428400000001
  $ws.Range("N39").Value = -7227.428400000001
  # Row 44: No More Dumpster Diving | Knight's Bread
  $ws.Range("H44").Value = 2767.3333
  $ws.Range("I44").Value = 0
  $ws.Range("K44").Value = 0
  $ws.Range("M44").ClearContents()
  # Row 50: Moving Up in the World | Rolanberry Cheese
  $ws.Range("H50").Value = 118.75
  $ws.Range("I50").Value = 91.666664
  $ws.Range("J50").Value = 200
  $ws.Range("K50").Value = 274.999992
  $ws.Range("L50").Value = 600
  $ws.Range("M50").Value = 206.000008
  $ws.Range("N50").Value = -1562
  # Row 53: Rolanberry Fields Forever | Rolanberry Cheese
  $ws.Range("H53").Value = 118.75
  $ws.Range("I53").Value = 91.666664
  $ws.Range("J53").Value = 200
  $ws.Range("K53").Value = 274.999992
  $ws.Range("L53").Value = 600
  $ws.Range("M53").Value = 206.000008
  $ws.Range("N53").Value = -1562
  # Row 55: Pagan Pastries | Pastry Fish
  $ws.Range("H55").Value = 3498.75
  $ws.Range("J55").Value = 3498.75
  $ws.Range("L55").Value = 10496.25
  $ws.Range("N55").Value = -10850.25
  # Row 62: Little Orphan Candy | Fig Bavarois
  $ws.Range("H62").Value = 5138
  $ws.Range("J62").Value = 5138
  $ws.Range("L62").Value = 15414
  $ws.Range("N62").Value = -16786
  # Row 65: Confections of Confession (L) | Fig Bavarois
  $ws.Range("H65").Value = 5138
  $ws.Range("J65").Value = 5138
  $ws.Range("L65").Value = 46242
  $ws.Range("N65").Value = -53106
  # Row 96: Hunger Is No Game | Popoto Soba
  $ws.Range("H96").Value = 7180.5713
  $ws.Range("I96").Value = 0
  $ws.Range("J96").Value = 7180.5713
  $ws.Range("K96").Value = 0
  $ws.Range("L96").Value = 21541.7139
  $ws.Range("M96").ClearContents()
  $ws.Range("N96").Value = -25659.7139
  # Row 107: Slippery Service | Frantoio Oil
  $ws.Range("H107").Value = 6286.2354
  $ws.Range("J107").Value = 8742.416999999999
  $ws.Range("L107").Value = 26227.251
  $ws.Range("N107").Value = -30067.251
  # Row 123: Topping Up the Pot | Zurek
  $ws.Range("H123").Value = 2989.8333
  $ws.Range("I123").Value = 2895
  $ws.Range("J123").Value = 3008.8
  $ws.Range("K123").Value = 8685
  $ws.Range("L123").Value = 9026.400000000001
  $ws.Range("M123").Value = -6235
  $ws.Range("N123").Value = -13926.4
  # Row 130: Blast from the Pasta | The Noodles of Elpis
  $ws.Range("H130").Value = 0
  $ws.Range("J130").Value = 0
  $ws.Range("L130").Value = 0
  $ws.Range("N130").ClearContents()
  # Row 131: The Mountain Steeped | Tsai tou Vounou
  $ws.Range("H131").Value = 22225360
  $ws.Range("I131").Value = 90909320
  $ws.Range("J131").Value = 4077.7646
  $ws.Range("K131").Value = 272727960
  $ws.Range("L131").Value = 12233.2938
  $ws.Range("M131").Value = -272722920
  $ws.Range("N131").Value = -22313.2938
  # Row 136: Simple Is Hardest | Spaghetti al Olio e Peperoncino
  $ws.Range("H136").Value = 1825.25
  $ws.Range("I136").Value = 883.9
  $ws.Range("K136").Value = 2651.7
  $ws.Range("M136").Value = 2448.3
  # Row 139: Najoothie | Wild Banana Blend
  $ws.Range("H139").Value = 1900.3
  $ws.Range("I139").Value = 2034.3334
  $ws.Range("J139").Value = 1699.25
  $ws.Range("K139").Value = 6103.0002
  $ws.Range("L139").Value = 5097.75
  $ws.Range("M139").Value = -963.0002000000004
  $ws.Range("N139").Value = -15377.75

$ws = $wb.Worksheets.Item("GSM")
  # Row 2: Copper and Robbers | Copper Ingot
  $ws.Range("H2").Value = 303
  $ws.Range("I2").Value = 303
  $ws.Range("J2").Value = 0
  $ws.Range("K2").Value = 303
  $ws.Range("L2").Value = 0
  $ws.Range("M2").Value = -190
  $ws.Range("N2").ClearContents()
  # Row 18: Gorgeous Gorget | Brass Gorget
  $ws.Range("H18").Value = 3006
  $ws.Range("J18").Value = 3006
  $ws.Range("L18").Value = 3006
  $ws.Range("N18").Value = -3592
  # Row 80: Needs More Prayerbell | Hardsilver Ingot
  $ws.Range("H80").Value = 3188
  $ws.Range("I80").Value = 1754.2858
  $ws.Range("K80").Value = 1754.2858
  $ws.Range("M80").Value = -756.2858000000001
  # Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
  $ws.Range("H83").Value = 3188
  $ws.Range("I83").Value = 1754.2858
  $ws.Range("K83").Value = 8771.429
  $ws.Range("M83").Value = -3779.429

$ws = $wb.Worksheets.Item("LTW")
  # Row 134: Freezing Fingers | Crocodileskin Fingerless Gloves of Striking
  $ws.Range("H134").Value = 30160
  $ws.Range("J134").Value = 30160
  $ws.Range("L134").Value = 30160
  $ws.Range("N134").Value = -40300

$ws = $wb.Worksheets.Item("WVR")
  # Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
  $ws.Range("H81").Value = 1125.25
  $ws.Range("I81").Value = 999.5
  $ws.Range("J81").Value = 1251
  $ws.Range("K81").Value = 1999
  $ws.Range("L81").Value = 2502
  $ws.Range("M81").Value = -938
  $ws.Range("N81").Value = -4624
  # Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
  $ws.Range("H84").Value = 1125.25
  $ws.Range("I84").Value = 999.5
  $ws.Range("J84").Value = 1251
  $ws.Range("K84").Value = 9995
  $ws.Range("L84").Value = 12510
  $ws.Range("M84").Value = -4691
  $ws.Range("N84").Value = -23118
  # Row 136: Weaving the Envelope | Sarcenet Cloth
  $ws.Range("H136").Value = 574.75
  $ws.Range("I136").Value = 266.25
  $ws.Range("K136").Value = 798.75
  $ws.Range("M136").Value = 1751.25

